$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B23: convert from text "3" to a real number 3
$ws.Range("B23").Value = 3

# Add new row 24 with annotation data
$ws.Range("A24").Value = "Sunsi Wu"
# B24 must stay a text string "2" (not numeric), so force text format first,
# then restore the default "General" format afterward to avoid leaving a
# stray number-format style on the cell.
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "2"
$ws.Range("B24").ClearFormats()
$ws.Range("C24").Value = "does not"
$ws.Range("D24").Value = "CRT"
$ws.Range("E24").Value = "MET"
$ws.Range("F24").Value = "b49eb73e-9ff0-45de-a177-7d78dc315c92"
$ws.Range("G24").Value = "2rHk2kZ5knTJ6_annotated.xlsx"
$ws.Range("H24").Value = "Cons: - it does not seem to involve any learning, it clearly does not fit at ICLR."
